$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the leading "#----...----" banner line (first paragraph):
#    delete its run's text entirely, leaving an empty paragraph.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Text = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$p1Text.Text = ""

# ------------------------------------------------------------------
# 2) "# Name : Susmitha " -> " Name : " + "Susmitha " (two runs,
#    same formatting) by first stripping the leading "#", then
#    forcing a run split right after " Name : ".
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$nameStart = $p2.Range.Start
$hashRange = $d.Range($nameStart, $nameStart + 1)
$hashRange.Text = ""

# " Name : Susmitha " now occupies [$nameStart, $nameStart+17)
# Force a genuine run boundary at $nameStart+8 (after " Name : ")
# by re-assigning FormattedText onto itself.
$splitRange = $d.Range($nameStart, $nameStart + 8)
$splitRange.FormattedText = $splitRange.FormattedText

# ------------------------------------------------------------------
# 3) "# Date : 11-03-2023" -> "Date : 11-03-2023"
# ------------------------------------------------------------------
$d.Content.Find.Execute("# Date : 11-03-2023", $true, $false, $false, $false, $false, $true, 1, $false, "Date : 11-03-2023", 2)

# ------------------------------------------------------------------
# 4) "# Subject : INTRO TO PROG USING SCRIPTING CS-504-D" ->
#    " Subject : INTRO TO PROG USING SCRIPTING CS-504-D"
# ------------------------------------------------------------------
$d.Content.Find.Execute("# Subject : INTRO TO PROG USING SCRIPTING CS-504-D", $true, $false, $false, $false, $false, $true, 1, $false, " Subject : INTRO TO PROG USING SCRIPTING CS-504-D", 2)

# ------------------------------------------------------------------
# 5) "# Title : Algorithm Assignment" -> "Title : Algorithm Assignment"
# ------------------------------------------------------------------
$d.Content.Find.Execute("# Title : Algorithm Assignment", $true, $false, $false, $false, $false, $true, 1, $false, "Title : Algorithm Assignment", 2)

# ------------------------------------------------------------------
# 6) "#----...----" (67 dashes) -> "----...----" (66 dashes)
# ------------------------------------------------------------------
$d.Content.Find.Execute("#-------------------------------------------------------------------", $true, $false, $false, $false, $false, $true, 1, $false, "------------------------------------------------------------------", 2)

# ------------------------------------------------------------------
# 7) styles.xml: add <w:semiHidden/> to the DefaultParagraphFont
#    character style (before <w:unhideWhenUsed/>).
# ------------------------------------------------------------------
$styles = $d.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.SemiHidden = $true

Write-Output "done"
